$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1352, 1).Value = "U46_01"
$ws.Cells.Item(1353, 1).Value = "U46_02"
$ws.Cells.Item(1354, 1).Value = "U46_03"
$ws.Cells.Item(1355, 1).Value = "U46_04"
$ws.Cells.Item(1356, 1).Value = "U46_05"
$ws.Cells.Item(1357, 1).Value = "U46_06"
$ws.Cells.Item(1358, 1).Value = "U46_07"
$ws.Cells.Item(1359, 1).Value = "U46_08"
$ws.Cells.Item(1360, 1).Value = "U46_09"
$ws.Cells.Item(1361, 1).Value = "U46_10"
$ws.Cells.Item(1362, 1).Value = "U46_11"
$ws.Cells.Item(1363, 1).Value = "U46_12"
$ws.Cells.Item(1364, 1).Value = "U46_13"
$ws.Cells.Item(1365, 1).Value = "U46_14"
$ws.Cells.Item(1366, 1).Value = "U46_15"
$ws.Cells.Item(1367, 1).Value = "U46_16"
$ws.Cells.Item(1368, 1).Value = "U46_17"
$ws.Cells.Item(1369, 1).Value = "U46_18"
$ws.Cells.Item(1370, 1).Value = "U46_19"
$ws.Cells.Item(1371, 1).Value = "U46_20"
$ws.Cells.Item(1372, 1).Value = "U46_21"
$ws.Cells.Item(1373, 1).Value = "U46_22"
$ws.Cells.Item(1374, 1).Value = "U46_23"
$ws.Cells.Item(1375, 1).Value = "U46_24"
$ws.Cells.Item(1376, 1).Value = "U46_25"
$ws.Cells.Item(1377, 1).Value = "U46_26"
$ws.Cells.Item(1378, 1).Value = "U46_27"
$ws.Cells.Item(1379, 1).Value = "U46_28"
$ws.Cells.Item(1380, 1).Value = "U46_29"
$ws.Cells.Item(1381, 1).Value = "U46_30"
$ws.Cells.Item(1352, 2).Value = 46
$ws.Cells.Item(1353, 2).Value = 46
$ws.Cells.Item(1354, 2).Value = 46
$ws.Cells.Item(1355, 2).Value = 46
$ws.Cells.Item(1356, 2).Value = 46
$ws.Cells.Item(1357, 2).Value = 46
$ws.Cells.Item(1358, 2).Value = 46
$ws.Cells.Item(1359, 2).Value = 46
$ws.Cells.Item(1360, 2).Value = 46
$ws.Cells.Item(1361, 2).Value = 46
$ws.Cells.Item(1362, 2).Value = 46
$ws.Cells.Item(1363, 2).Value = 46
$ws.Cells.Item(1364, 2).Value = 46
$ws.Cells.Item(1365, 2).Value = 46
$ws.Cells.Item(1366, 2).Value = 46
$ws.Cells.Item(1367, 2).Value = 46
$ws.Cells.Item(1368, 2).Value = 46
$ws.Cells.Item(1369, 2).Value = 46
$ws.Cells.Item(1370, 2).Value = 46
$ws.Cells.Item(1371, 2).Value = 46
$ws.Cells.Item(1372, 2).Value = 46
$ws.Cells.Item(1373, 2).Value = 46
$ws.Cells.Item(1374, 2).Value = 46
$ws.Cells.Item(1375, 2).Value = 46
$ws.Cells.Item(1376, 2).Value = 46
$ws.Cells.Item(1377, 2).Value = 46
$ws.Cells.Item(1378, 2).Value = 46
$ws.Cells.Item(1379, 2).Value = 46
$ws.Cells.Item(1380, 2).Value = 46
$ws.Cells.Item(1381, 2).Value = 46
$ws.Cells.Item(1352, 3).Value = "Ngôn ngữ"
$ws.Cells.Item(1352, 4).Value = "Language"
$ws.Cells.Item(1352, 5).Value = "English is a foreign language"
$ws.Cells.Item(1352, 6).Value = "foreign language "
$ws.Cells.Item(1352, 7).Value = "N"
$ws.Cells.Item(1353, 3).Value = "Thuộc nước ngoài"
$ws.Cells.Item(1353, 4).Value = "Foreign"
$ws.Cells.Item(1353, 5).Value = "We haven't been to foreign countries"
$ws.Cells.Item(1353, 6).Value = "a foreign country / nước ngoài"
$ws.Cells.Item(1353, 7).Value = "Adj"
$ws.Cells.Item(1354, 3).Value = "Hành trình, chuyến đi"
$ws.Cells.Item(1354, 4).Value = "Journey"
$ws.Cells.Item(1354, 5).Value = "Take me on a journey with you"
$ws.Cells.Item(1354, 6).Value = "on a journey"
$ws.Cells.Item(1354, 7).Value = "N"
$ws.Cells.Item(1355, 3).Value = "Người mới bắt đầu"
$ws.Cells.Item(1355, 4).Value = "Beginner"
$ws.Cells.Item(1355, 5).Value = "This level is easy for a beginner"
$ws.Cells.Item(1355, 6).Value = "for a beginner"
$ws.Cells.Item(1355, 7).Value = "N"
$ws.Cells.Item(1356, 3).Value = "Từ, lời nói"
$ws.Cells.Item(1356, 4).Value = "Word"
$ws.Cells.Item(1356, 5).Value = "My leader has the last word on team decisions (quyết định)"
$ws.Cells.Item(1356, 6).Value = "have the last word on something / đưa lời nói cuối cùng về việc gì đó"
$ws.Cells.Item(1356, 7).Value = "N"
$ws.Cells.Item(1357, 3).Value = "Lưu loát"
$ws.Cells.Item(1357, 4).Value = "Fluent"
$ws.Cells.Item(1357, 5).Value = "she is fluent in French"
$ws.Cells.Item(1357, 6).Value = "fluent in something"
$ws.Cells.Item(1357, 7).Value = "Adj"
$ws.Cells.Item(1358, 3).Value = "Bản địa"
$ws.Cells.Item(1358, 4).Value = "Native"
$ws.Cells.Item(1358, 6).Value = "Native language / tiếng mẹ đẻ"
$ws.Cells.Item(1358, 5).Value = "My native language is Vietnamese"
$ws.Cells.Item(1358, 7).Value = "Adj"
$ws.Cells.Item(1359, 3).Value = "người nói"
$ws.Cells.Item(1359, 4).Value = "Speaker"
$ws.Cells.Item(1359, 5).Value = "She was a brilliant public speaker"
$ws.Cells.Item(1359, 6).Value = "a public speaker / một diễn giả"
$ws.Cells.Item(1359, 7).Value = "N"
$ws.Cells.Item(1360, 3).Value = "Dài hạn"
$ws.Cells.Item(1360, 4).Value = "long-term"
$ws.Cells.Item(1360, 5).Value = "A long-term goal often takes years to achieve"
$ws.Cells.Item(1360, 6).Value = "long-term goal"
$ws.Cells.Item(1360, 7).Value = "Adj"
$ws.Cells.Item(1361, 3).Value = "Khái niệm, ý tưởng chủ đạo"
$ws.Cells.Item(1361, 4).Value = "Concept"
$ws.Cells.Item(1361, 5).Value = "The concept of this event is experimental (Đang được thử nghiệm)"
$ws.Cells.Item(1361, 6).Value = "Concept of something / ý tưởng chủ đạo của something"
$ws.Cells.Item(1361, 7).Value = "N"
$ws.Cells.Item(1362, 3).Value = "Ghi lại"
$ws.Cells.Item(1362, 4).Value = "Note"
$ws.Cells.Item(1362, 5).Value = "I often note down new words in my handbook"
$ws.Cells.Item(1362, 6).Value = "note down something"
$ws.Cells.Item(1362, 7).Value = "V"
$ws.Cells.Item(1363, 3).Value = "Đánh vần"
$ws.Cells.Item(1363, 4).Value = "Spell"
$ws.Cells.Item(1363, 5).Value = "Please spell your name for me"
$ws.Cells.Item(1363, 6).Value = "spell something / đánh vần một cái gì đó"
$ws.Cells.Item(1363, 7).Value = "V"
$ws.Cells.Item(1364, 3).Value = "Ngữ pháp"
$ws.Cells.Item(1364, 4).Value = "Grammar"
$ws.Cells.Item(1364, 5).Value = "English grammar is easy to learn"
$ws.Cells.Item(1364, 6).Value = "English grammar"
$ws.Cells.Item(1364, 7).Value = "N"
$ws.Cells.Item(1365, 3).Value = "Kỷ thuật"
$ws.Cells.Item(1365, 4).Value = "Technique"
$ws.Cells.Item(1365, 5).Value = "You must learn the technique for baking."
$ws.Cells.Item(1365, 6).Value = "technique for doing something / kỷ thuật làm một cái gì đó"
$ws.Cells.Item(1365, 7).Value = "N"
$ws.Cells.Item(1366, 3).Value = "Thường xuyên"
$ws.Cells.Item(1366, 4).Value = "Regular"
$ws.Cells.Item(1366, 5).Value = "I am the restaurant's regular customer"
$ws.Cells.Item(1366, 6).Value = "a regular customer / khách quen"
$ws.Cells.Item(1366, 7).Value = "Adj"
$ws.Cells.Item(1367, 3).Value = "Dịch"
$ws.Cells.Item(1367, 4).Value = "Translate"
$ws.Cells.Item(1367, 5).Value = "Can you translate english into vietnamese?"
$ws.Cells.Item(1367, 6).Value = "Translate something into something / dịch một cái gì thành một cái gì"
$ws.Cells.Item(1367, 7).Value = "V"
$ws.Cells.Item(1368, 3).Value = "Bối cảnh"
$ws.Cells.Item(1368, 4).Value = "Context"
$ws.Cells.Item(1368, 6).Value = "in or within the context of something / trong bối cảnh gì đó"
$ws.Cells.Item(1368, 5).Value = "In the context of funerals(đám tang), wear dark clothing."
$ws.Cells.Item(1368, 7).Value = "N"
$ws.Cells.Item(1369, 3).Value = "Định nghĩa"
$ws.Cells.Item(1369, 4).Value = "Define"
$ws.Cells.Item(1369, 5).Value = "Please define this work for me"
$ws.Cells.Item(1369, 6).Value = "define a word / định nghĩa một từ"
$ws.Cells.Item(1369, 7).Value = "V"
$ws.Cells.Item(1370, 3).Value = "Từ điển"
$ws.Cells.Item(1370, 4).Value = "Dictionary"
$ws.Cells.Item(1370, 5).Value = "You can look it up in the dictionary"
$ws.Cells.Item(1370, 6).Value = "look it up in the dictionary / tra cứu trong từ điển"
$ws.Cells.Item(1370, 7).Value = "N"
$ws.Cells.Item(1371, 3).Value = "Động từ"
$ws.Cells.Item(1371, 4).Value = "Verb"
$ws.Cells.Item(1371, 5).Value = "`"To be`" is an irregular verb "
$ws.Cells.Item(1371, 6).Value = "irregular verb / động từ bất quy tắc"
$ws.Cells.Item(1371, 7).Value = "N"
$ws.Cells.Item(1372, 3).Value = "Tương tác"
$ws.Cells.Item(1372, 4).Value = "Interact"
$ws.Cells.Item(1372, 5).Value = "I love interacting with people"
$ws.Cells.Item(1372, 6).Value = "interact with somebody / tương tác với ai đó"
$ws.Cells.Item(1372, 7).Value = "V"
$ws.Cells.Item(1373, 3).Value = "Cụm từ"
$ws.Cells.Item(1373, 4).Value = "Phrase"
$ws.Cells.Item(1373, 5).Value = "Catchy phrase help me study better"
$ws.Cells.Item(1373, 6).Value = "catchy phrase / cụm từ dễ nhớ"
$ws.Cells.Item(1373, 7).Value = "N"
$ws.Cells.Item(1374, 3).Value = "Giảm tốc độ"
$ws.Cells.Item(1374, 4).Value = "Slow down"
$ws.Cells.Item(1374, 5).Value = "slow the car down"
$ws.Cells.Item(1374, 6).Value = "slow somebody of something down "
$ws.Cells.Item(1374, 7).Value = "V"
$ws.Cells.Item(1375, 3).Value = "Căn bản, cơ bản"
$ws.Cells.Item(1375, 4).Value = "Basic"
$ws.Cells.Item(1375, 5).Value = "Lions have basic instincts in hunting"
$ws.Cells.Item(1375, 6).Value = "basic instinct / bản năng căn bản"
$ws.Cells.Item(1375, 7).Value = "Adj"
$ws.Cells.Item(1376, 3).Value = "Quên"
$ws.Cells.Item(1376, 4).Value = "Forget"
$ws.Cells.Item(1376, 5).Value = "I forgot about you"
$ws.Cells.Item(1376, 6).Value = "forget about something or someone / quên đi một cái gì đó hoặc ai đó"
$ws.Cells.Item(1376, 7).Value = "V"
$ws.Cells.Item(1377, 4).Value = "Sound"
$ws.Cells.Item(1377, 3).Value = "Nghe"
$ws.Cells.Item(1377, 5).Value = "Does my accent sound right?"
$ws.Cells.Item(1377, 6).Value = "sound right / nghe đúng, nghe chuẩn"
$ws.Cells.Item(1377, 7).Value = "V"
$ws.Cells.Item(1378, 3).Value = "Trò chuyện"
$ws.Cells.Item(1378, 4).Value = "Chat"
$ws.Cells.Item(1378, 5).Value = "She is chatting with someone at the bar"
$ws.Cells.Item(1378, 6).Value = "chat with somebody"
$ws.Cells.Item(1378, 7).Value = "V"
$ws.Cells.Item(1379, 3).Value = "Câu"
$ws.Cells.Item(1379, 4).Value = "Sentence"
$ws.Cells.Item(1379, 5).Value = "English only has a few basic sentence structures"
$ws.Cells.Item(1379, 6).Value = "sentence structure / cấu trúc câu"
$ws.Cells.Item(1379, 7).Value = "N"
$ws.Cells.Item(1380, 3).Value = "Tự tin"
$ws.Cells.Item(1380, 4).Value = "Confident"
$ws.Cells.Item(1380, 5).Value = "Are you confident about your performance?"
$ws.Cells.Item(1380, 6).Value = "confident about something"
$ws.Cells.Item(1380, 7).Value = "Adj"
$ws.Cells.Item(1381, 3).Value = "Ý nghĩa"
$ws.Cells.Item(1381, 4).Value = "Meaning"
$ws.Cells.Item(1381, 5).Value = "What is the meaning of life ="
$ws.Cells.Item(1381, 6).Value = "meaning of something"
$ws.Cells.Item(1381, 7).Value = "N"

$ws.Range("F1366").Select() | Out-Null
